$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '51.778.30'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  +4.82%  '

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '2.770.46'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  +5.37%  '

$ws.Range("E4").Value = '  +0.03%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '116.99'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +4.69%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '333.12'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +2.66%  '

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.539'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  +2.65%  '

$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("E9").Value = '  +6.19%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '42.24'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +6.89%  '

$ws.Range("E11").Value = '  +6.04%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '20.30'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +2.48%  '

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '0.130'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +2.06%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '7.65'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +4.72%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '3.205.71'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +5.38%  '

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '2.775.62'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  +5.88%  '

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '0.891'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  +5.30%  '

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '51.688.87'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +4.80%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '3.26'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +10.90%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '13.58'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +5.67%  '

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '6.86'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +2.67%  '

$ws.Range("E22").Value = '  +3.34%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '279.70'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  +3.73%  '

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '70.04'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  +1.77%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '2.70'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +6.86%  '

$ws.Range("E26").Value = '  +2.55%  '

$ws.Range("E27").Value = '  +0.04%  '

$ws.Range("E28").Value = '  -0.59%  '

$ws.Range("E29").Value = '  +0.77%  '

$ws.Range("E30").Value = '  +3.80%  '

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '35.14'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +1.53%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '50.12'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +1.41%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '5.59'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +2.07%  '

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '0.0822'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +1.00%  '

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '19.22'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  +1.89%  '

$ws.Range("E36").Value = '  -0.09%  '

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '5.02'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  +2.93%  '

$ws.Range("E38").Value = '  +3.03%  '

$ws.Range("E39").Value = '  +5.02%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.0358'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  +9.85%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '127.19'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -1.36%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '23.34'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +6.35%  '

$ws.Range("E43").Value = '  +8.21%  '

$ws.Range("E44").Value = '  +3.06%  '

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '2.46'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +16.63%  '

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '2.088.22'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  +1.66%  '

$ws.Range("E47").Value = '  +3.91%  '

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '2.24'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  +4.19%  '

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '5.53'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  +6.76%  '

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '60.69'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  +2.78%  '

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '8.84'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -0.70%  '
